$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 225, pushing the existing rows 225-316 down to 226-317.
$ws.Rows.Item(225).Insert()

# Populate the newly inserted row 225 with the new price record.
$ws.Cells.Item(225, 1).Value = 7
$ws.Cells.Item(225, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(225, 3).Value = "Ñuble"
$ws.Cells.Item(225, 4).Value = 44510
$ws.Cells.Item(225, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(225, 5).Value = 16
$ws.Cells.Item(225, 6).Value = 100112020
$ws.Cells.Item(225, 7).Value = "Tomate"
$ws.Cells.Item(225, 8).Value = "Larga vida"
$ws.Cells.Item(225, 9).Value = "Primera"
$ws.Cells.Item(225, 10).Value = 240
$ws.Cells.Item(225, 11).Value = 7500
$ws.Cells.Item(225, 12).Value = 8000
$ws.Cells.Item(225, 13).Value = 7750
$ws.Cells.Item(225, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(225, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(225, 16).Value = 775
$ws.Cells.Item(225, 17).Value = 10
$ws.Cells.Item(225, 18).Value = "Hortaliza"
